# "change keyboard press to mouse press"
# Updates the Settings sheet (resolution ratio, repeat count, pre-first-side
# delay, the assist-servant equipment note, and a new emulator-window-name
# row), then rebalances which skills are flagged "use" on the Side1/Side2/
# Side3 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("B2").Value = 1.25
$settings.Range("B3").Value = 7
$settings.Range("B7").Value = 15
$settings.Range("B6").Value = "bondage"

$settings.Range("A9").Copy()
$settings.Range("A10").PasteSpecial(-4122)
$settings.Range("A10").Value = "模拟器窗口名称（默认为网易mumu模拟器）："
$settings.Range("B10").Value = "命运-冠位指定 - MuMu模拟器"

# ---------------------------------------------------------------------
# Side1 sheet - clear a few previously-enabled skills
# ---------------------------------------------------------------------
$side1 = $wb.Worksheets.Item("Side1")
$side1.Range("B7").Value = $null
$side1.Range("C7").Value = $null
$side1.Range("B8").Value = $null
$side1.Range("C8").Value = $null
$side1.Range("B16").Value = $null
$side1.Range("C16").Value = $null

# ---------------------------------------------------------------------
# Side2 sheet - enable several skills
# ---------------------------------------------------------------------
$side2 = $wb.Worksheets.Item("Side2")
$side2.Range("B5").Value = 1
$side2.Range("C5").Value = 0
$side2.Range("B8").Value = 1
$side2.Range("C8").Value = 0
$side2.Range("B9").Value = 1
$side2.Range("C9").Value = 1
$side2.Range("B10").Value = 1
$side2.Range("C10").Value = 0
$side2.Range("B11").Value = 1
$side2.Range("C11").Value = 0

# ---------------------------------------------------------------------
# Side3 sheet - enable/disable skills and restore page setup
# ---------------------------------------------------------------------
$side3 = $wb.Worksheets.Item("Side3")
$side3.Range("B6").Value = 1
$side3.Range("C6").Value = 0
$side3.Range("B13").Value = 1
$side3.Range("C13").Value = 2
$side3.Range("B16").Value = $null
$side3.Range("C16").Value = $null
$side3.Range("B17").Value = 1
$side3.Range("C17").Value = 0

$side3.PageSetup.PaperSize = 9
$side3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selection bookkeeping (matches the cursor positions left behind by the
# author's edit session)
# ---------------------------------------------------------------------
$side1.Range("C27").Select() | Out-Null
$side2.Range("C16").Select() | Out-Null
$side3.Range("A16").Select() | Out-Null
$settings.Range("A9").Select() | Out-Null
